# SMP_3PP.xlsx edit: update assumed particle/void size input (AP3:AP18)
# from 200000 to 136735.30180419001, add a new summary row 23 with
# W23 = MIN(W2:W22), and move the active selection to Y3 (scrolled so
# column R is visible at the left).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update AP3:AP18 input values (all rows shared the same constant) ---
$newValue = 136735.30180419001
for ($row = 3; $row -le 18; $row++) {
    $ws.Cells.Item($row, 42).Value = $newValue   # column AP = 42
}

# --- Add new row 23: W23 = MIN(W2:W22) ---
$ws.Range("W23").Formula = "=MIN(W2:W22)"

# --- Recalculate so all dependent cells (AQ, AU, AV, AW, AX, ...) refresh ---
$excel.Calculate()

# --- Move the view / selection like in the committed workbook ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 18   # column R
$win.ScrollRow = 1
$ws.Range("Y3").Select()
